# dc_pod_v010 - "missing factors fixed"
# Adds 4 new data rows (142-145) to Tier_I, fills in the previously-missing
# seasonally_adj ("n") marker for rows 131-141, and adds a new "Sheet1"
# worksheet holding the "global groups" legend.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Backfill the missing "seasonally_adj" column (G) for rows 131-141.
#    These rows already exist but were missing the G cell entirely.
# ---------------------------------------------------------------------
for ($r = 131; $r -le 141; $r++) {
  $ws1.Range("G$r").Value = "n"
}

# ---------------------------------------------------------------------
# 2. Append the four new source rows (142-145).
#    Column order matches the authoring sequence so that new shared
#    strings land in the same order as the source edit:
#      A142, A143 (tickers) -> then the Sheet1 legend -> then D/B pairs.
# ---------------------------------------------------------------------
$ws1.Range("A142").Value = "IUES.L"
$ws1.Range("A143").Value = "IUIT.L"

# ---------------------------------------------------------------------
# 3. New worksheet "Sheet1" (inserted right after Tier_I) holding the
#    "global groups" legend used elsewhere in the model.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("B4").Value  = "global groups:"
$ws2.Range("C5").Value  = "> macro: economic activity"
$ws2.Range("C6").Value  = "> macro: demand"
$ws2.Range("C7").Value  = "> macro: revenues"
$ws2.Range("C8").Value  = "> macro: liquidity impulse"
$ws2.Range("C9").Value  = "> commodities: prices"
$ws2.Range("C10").Value = "> capital: yields in different classes"

# ---------------------------------------------------------------------
# 4. Finish filling in the new Tier_I rows, column by column, in the
#    same order the strings were originally authored.
# ---------------------------------------------------------------------
$ws1.Range("D142").Value = "https://finance.yahoo.com/quote/IUES.L/history/"
$ws1.Range("B142").Value = "iShares S&P 500 Energy Sector UCITS ETF USD (Acc)"

$ws1.Range("D143").Value = "https://finance.yahoo.com/quote/IUIT.L/history/"
$ws1.Range("B143").Value = "iShares S&P 500 Information Technology Sector UCITS ETF USD (Acc)"

$ws1.Range("D144").Value = "https://www.etf.com/etfanalytics/etf-fund-flows-tool-result?tickers=IVV%2C&startDate=2007-07-02&endDate=2024-08-06&frequency=MONTHLY"
$ws1.Range("B144").Value = "Calculated manually as flows(VOO + IVV + SPY + VTI + QQQ)"
$ws1.Range("A144").Value = "SNPF"
$ws1.Range("F144").Value = "W-MON"

$ws1.Range("A145").Value = "DTWEXBGS"
$ws1.Range("B145").Value = "Nominal Broad U.S. Dollar Index"
$ws1.Range("D145").Value = "https://fred.stlouisfed.org/series/DTWEXBGS"

# ---------------------------------------------------------------------
# 5. Remaining columns for the new rows: value_type (C), reader (E) and
#    ts_frequency (F) for rows 142/143/145 - all reuse existing values.
# ---------------------------------------------------------------------
$ws1.Range("C142").Value = "CONTINUOUS"
$ws1.Range("E142").Value = "csv"
$ws1.Range("F142").Value = "MS"
$ws1.Range("H142").Value = 0

$ws1.Range("C143").Value = "CONTINUOUS"
$ws1.Range("E143").Value = "csv"
$ws1.Range("F143").Value = "MS"
$ws1.Range("H143").Value = 0

$ws1.Range("C144").Value = "CONTINUOUS"
$ws1.Range("E144").Value = "csv"
$ws1.Range("H144").Value = 0

$ws1.Range("C145").Value = "CONTINUOUS"
$ws1.Range("E145").Value = "csv"
$ws1.Range("F145").Value = "B"
$ws1.Range("H145").Value = 0

# seasonally_adj ("n", red highlight) for the new rows too.
for ($r = 142; $r -le 145; $r++) {
  $ws1.Range("G$r").Value = "n"
}

# ---------------------------------------------------------------------
# 6. Copy the red "n" formatting (style used throughout column G) onto
#    every G cell we just touched, matching style index 4 in the file.
# ---------------------------------------------------------------------
$ws1.Range("G50").Copy()
$ws1.Range("G131:G145").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 7. New column J width (only column J differs from the rest of F:R).
# ---------------------------------------------------------------------
$ws1.Columns.Item(10).ColumnWidth = 19

# ---------------------------------------------------------------------
# 8. Selections / active sheet, matching the saved view state.
# ---------------------------------------------------------------------
$ws2.Range("C11").Select()
$ws1.Range("J145").Select()
$ws1.Activate()
